$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "https://stackoverflow.com/questions/43637687/python-automation?r=SearchResults"
$ws.Range("A4").Value = "https://stackoverflow.com/questions/63717551/how-to-use-python-automation-in-flutter-app?r=SearchResults"
$ws.Range("A5").Value = "https://stackoverflow.com/questions/51451268/error-while-exeutiing-selenium-python-automation-script?r=SearchResults"
$ws.Range("A6").Value = "https://stackoverflow.com/questions/60994919/deploying-a-python-automation-script-in-the-cloud?r=SearchResults"
$ws.Range("A7").Value = "https://stackoverflow.com/questions/40208051/selenium-using-python-geckodriver-executable-needs-to-be-in-path?r=SearchResults"
$ws.Range("A8").Value = "https://stackoverflow.com/questions/37048354/python-automation-for-android?r=SearchResults"
$ws.Range("A9").Value = "https://stackoverflow.com/questions/71252278/calculating-server-throughput-for-selenium-python-automation-test?r=SearchResults"
$ws.Range("A10").Value = "https://stackoverflow.com/questions/77147730/python-automation-using-selenium?r=SearchResults"
$ws.Range("A11").Value = "https://stackoverflow.com/questions/40914325/python-automation-using-subprocess?r=SearchResults"

$wb.Save()
